$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target cells keep their original plain-text representation
# (exact trailing zeros, percent signs, etc.) instead of being
# auto-coerced to numbers/percentages by Excel's type inference.
$targetCells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "E19", "D20", "E20", "D21", "E21", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D47", "E47")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# New Price (column D) / Volume(1h) (column E) values scraped on
# Thu Jan 12 23:28:05 UTC 2023.
$ws.Range("D2").Value = "286.65"
$ws.Range("E2").Value = "2.02%"
$ws.Range("D3").Value = "28.60"
$ws.Range("E3").Value = "3.59%"
$ws.Range("E4").Value = "4.29%"
$ws.Range("D5").Value = "0.06660"
$ws.Range("E5").Value = "3.04%"
$ws.Range("D6").Value = "7.389"
$ws.Range("E6").Value = "4.47%"
$ws.Range("D7").Value = "3.407"
$ws.Range("E7").Value = "3.07%"
$ws.Range("D8").Value = "1.372"
$ws.Range("E8").Value = "7.25%"
$ws.Range("D9").Value = "0.9352"
$ws.Range("E9").Value = "2.91%"
$ws.Range("E10").Value = "2.31%"
$ws.Range("D11").Value = "0.06475"
$ws.Range("E11").Value = "1.98%"
$ws.Range("D12").Value = "0.07559"
$ws.Range("E12").Value = "0.80%"
$ws.Range("D13").Value = "0.02937"
$ws.Range("E13").Value = "0.02%"
$ws.Range("D14").Value = "0.08989"
$ws.Range("E14").Value = "-0.15%"
$ws.Range("D15").Value = "0.001593"
$ws.Range("E15").Value = "-0.15%"
$ws.Range("D16").Value = "0.04499"
$ws.Range("E16").Value = "1.90%"
$ws.Range("D17").Value = "0.0006470"
$ws.Range("E17").Value = "0.95%"
$ws.Range("D18").Value = "0.006259"
$ws.Range("E18").Value = "3.96%"
$ws.Range("E19").Value = "-1.22%"
$ws.Range("D20").Value = "2.254"
$ws.Range("E20").Value = "0.91%"
$ws.Range("D21").Value = "0.3217"
$ws.Range("E21").Value = "2.30%"
$ws.Range("E22").Value = "-4.11%"
$ws.Range("D23").Value = "4.087"
$ws.Range("E23").Value = "4.82%"
$ws.Range("D24").Value = "0.1552"
$ws.Range("E24").Value = "3.33%"
$ws.Range("D25").Value = "0.001181"
$ws.Range("E25").Value = "0.45%"
$ws.Range("D26").Value = "0.004147"
$ws.Range("E26").Value = "-3.97%"
$ws.Range("D27").Value = "0.0001249"
$ws.Range("E27").Value = "5.98%"
$ws.Range("D28").Value = "0.0001617"
$ws.Range("E28").Value = "-2.34%"
$ws.Range("D40").Value = "0.04214"
$ws.Range("E40").Value = "2.97%"
$ws.Range("D41").Value = "0.006728"
$ws.Range("E41").Value = "-1.91%"
$ws.Range("D42").Value = "0.1247"
$ws.Range("E42").Value = "-11.12%"
$ws.Range("D43").Value = "0.002018"
$ws.Range("E43").Value = "-2.40%"
$ws.Range("D44").Value = "0.01217"
$ws.Range("E44").Value = "4.36%"
$ws.Range("D45").Value = "0.00005588"
$ws.Range("E45").Value = "1.05%"
$ws.Range("D47").Value = "0.01306"
$ws.Range("E47").Value = "-29.29%"
